$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'235.88"
$ws.Range("D2").Style = "Normal"
$ws.Range("G2").Value = "'9"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'22.37"
$ws.Range("D3").Style = "Normal"
$ws.Range("G3").Value = "'9"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'5.412"
$ws.Range("D4").Style = "Normal"
$ws.Range("G4").Value = "'9"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.05636"
$ws.Range("D5").Style = "Normal"
$ws.Range("G5").Value = "'9"
$ws.Range("G5").Style = "Normal"
$ws.Range("D6").Value = "'3.374"
$ws.Range("D6").Style = "Normal"
$ws.Range("G6").Value = "'9"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'6.474"
$ws.Range("D7").Style = "Normal"
$ws.Range("G7").Value = "'9"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'1.077"
$ws.Range("D8").Style = "Normal"
$ws.Range("G8").Value = "'9"
$ws.Range("G8").Style = "Normal"
$ws.Range("G9").Value = "'9"
$ws.Range("G9").Style = "Normal"
$ws.Range("B10").Value = "One"
$ws.Range("C10").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D10").Value = "'0.01176"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "9OneONEBestin24h"
$ws.Range("G10").Value = "'9"
$ws.Range("G10").Style = "Normal"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1399"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "10WazirXWRX"
$ws.Range("G11").Value = "'9"
$ws.Range("G11").Style = "Normal"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.07412"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "11MandalaExchangeTokenMDX"
$ws.Range("G12").Value = "'9"
$ws.Range("G12").Style = "Normal"
$ws.Range("B13").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C13").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D13").Value = "'0.03205"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "12LiechtensteinCryptoassetsExchangeLCX"
$ws.Range("G13").Value = "'9"
$ws.Range("G13").Style = "Normal"
$ws.Range("B14").Value = "BitrueCoin"
$ws.Range("C14").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D14").Value = "'0.02941"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "13BitrueCoinBTR"
$ws.Range("G14").Value = "'9"
$ws.Range("G14").Style = "Normal"
$ws.Range("B15").Value = "BitMartToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D15").Value = "'0.09257"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "14BitMartTokenBMX"
$ws.Range("G15").Value = "'9"
$ws.Range("G15").Style = "Normal"
$ws.Range("B16").Value = "BitForexToken"
$ws.Range("C16").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D16").Value = "'0.001669"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "15BitForexTokenBF"
$ws.Range("G16").Value = "'9"
$ws.Range("G16").Style = "Normal"
$ws.Range("B17").Value = "MCDex"
$ws.Range("C17").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D17").Value = "'3.263"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "16MCDexMCB"
$ws.Range("G17").Value = "'9"
$ws.Range("G17").Style = "Normal"
$ws.Range("B18").Value = "CoinExToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D18").Value = "'0.04751"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "17CoinExTokenCET"
$ws.Range("G18").Value = "'9"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'0.006207"
$ws.Range("D19").Style = "Normal"
$ws.Range("G19").Value = "'9"
$ws.Range("G19").Style = "Normal"
$ws.Range("D20").Value = "'0.005107"
$ws.Range("D20").Style = "Normal"
$ws.Range("G20").Value = "'9"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.001052"
$ws.Range("D21").Style = "Normal"
$ws.Range("G21").Value = "'9"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'0.0001501"
$ws.Range("D22").Style = "Normal"
$ws.Range("G22").Value = "'9"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'3.868"
$ws.Range("D23").Style = "Normal"
$ws.Range("G23").Value = "'9"
$ws.Range("G23").Style = "Normal"
$ws.Range("G24").Value = "'9"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.3324"
$ws.Range("D25").Style = "Normal"
$ws.Range("G25").Value = "'9"
$ws.Range("G25").Style = "Normal"
$ws.Range("G26").Value = "'9"
$ws.Range("G26").Style = "Normal"
$ws.Range("D27").Value = "'0.0004992"
$ws.Range("D27").Style = "Normal"
$ws.Range("G27").Value = "'9"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'9"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'9"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'9"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'9"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'9"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'9"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'9"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'9"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'9"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'9"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'9"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'9"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.04055"
$ws.Range("D40").Style = "Normal"
$ws.Range("G40").Value = "'9"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.006983"
$ws.Range("D41").Style = "Normal"
$ws.Range("G41").Value = "'9"
$ws.Range("G41").Style = "Normal"
$ws.Range("B42").Value = "CEJI"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D42").Value = "'0.003501"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "41CEJICEJI"
$ws.Range("G42").Value = "'9"
$ws.Range("G42").Style = "Normal"
$ws.Range("B43").Value = "BKEXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D43").Value = "'0.1039"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "42BKEXTokenBKK"
$ws.Range("G43").Value = "'9"
$ws.Range("G43").Style = "Normal"
$ws.Range("D44").Value = "'0.009305"
$ws.Range("D44").Style = "Normal"
$ws.Range("G44").Value = "'9"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005431"
$ws.Range("D45").Style = "Normal"
$ws.Range("G45").Value = "'9"
$ws.Range("G45").Style = "Normal"
$ws.Range("G46").Value = "'9"
$ws.Range("G46").Style = "Normal"
$ws.Range("D47").Value = "'0.6755"
$ws.Range("D47").Style = "Normal"
$ws.Range("G47").Value = "'9"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.03933"
$ws.Range("D48").Style = "Normal"
$ws.Range("G48").Value = "'9"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("G49").Value = "'9"
$ws.Range("G49").Style = "Normal"
$ws.Range("G50").Value = "'9"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'9"
$ws.Range("G51").Style = "Normal"
